$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -168.6
$ws.Range("B3").Value = -245.1
$ws.Range("C3").Value = -188.8
$ws.Range("C4").Value = -157.7
